# Add a "FRI headways and runtimes" worksheet, positioned right after the
# "WKDY headways and runtimes" sheet (and before "SAT headways and
# runtimes"). The new sheet is an exact duplicate of the WKDY sheet
# (same headers, data, column widths and styles) -- mirrors using Excel's
# "Move or Copy... > Create a copy" on the WKDY tab, then renaming the
# resulting copy to "FRI headways and runtimes".

$wb = $excel.ActiveWorkbook

$wkdy = $wb.Worksheets.Item("WKDY headways and runtimes")
$sat  = $wb.Worksheets.Item("SAT headways and runtimes")

# Copy WKDY, inserting the new sheet immediately before SAT (i.e. right
# after WKDY).
$wkdy.Copy($sat)

# The freshly-created copy becomes active and is named
# "WKDY headways and runtimes (2)"; rename it to the Friday tab.
$fri = $wb.ActiveSheet
$fri.Name = "FRI headways and runtimes"
